# Табель за Август - г.Алматы : add overtime row block for "ilyas"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Column width adjustments (B and D:AI)
# ColumnWidth set via COM lands ~0.8333333333333333 chars wider than
# the requested "display" width once saved to OOXML, so subtract that
# fixed offset to land exactly on the target width.
# -----------------------------------------------------------------
$offset = 0.8333333333333333
$widths = @{
    2  = 7
    4  = 22
    5  = 22
    6  = 22
    7  = 26
    8  = 26
    9  = 22
    10 = 20
    11 = 22
    12 = 22
    13 = 22
    14 = 26
    15 = 26
    16 = 22
    17 = 20
    18 = 22
    19 = 22
    20 = 22
    21 = 26
    22 = 26
    23 = 22
    24 = 20
    25 = 22
    26 = 22
    27 = 22
    28 = 26
    29 = 26
    30 = 22
    31 = 20
    32 = 22
    33 = 22
    34 = 22
    35 = 17
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - $offset
}

# -----------------------------------------------------------------
# 2. Row 2 header cells: append the Russian weekday name to each date,
#    and mark the overtime column with its unit.
# -----------------------------------------------------------------
$ws.Range("D2").Value = "2024-08-01 (Четверг)"
$ws.Range("E2").Value = "2024-08-02 (Пятница)"
$ws.Range("F2").Value = "2024-08-03 (Суббота)"
$ws.Range("G2").Value = "2024-08-04 (Воскресенье)"
$ws.Range("H2").Value = "2024-08-05 (Понедельник)"
$ws.Range("I2").Value = "2024-08-06 (Вторник)"
$ws.Range("J2").Value = "2024-08-07 (Среда)"
$ws.Range("K2").Value = "2024-08-08 (Четверг)"
$ws.Range("L2").Value = "2024-08-09 (Пятница)"
$ws.Range("M2").Value = "2024-08-10 (Суббота)"
$ws.Range("N2").Value = "2024-08-11 (Воскресенье)"
$ws.Range("O2").Value = "2024-08-12 (Понедельник)"
$ws.Range("P2").Value = "2024-08-13 (Вторник)"
$ws.Range("Q2").Value = "2024-08-14 (Среда)"
$ws.Range("R2").Value = "2024-08-15 (Четверг)"
$ws.Range("S2").Value = "2024-08-16 (Пятница)"
$ws.Range("T2").Value = "2024-08-17 (Суббота)"
$ws.Range("U2").Value = "2024-08-18 (Воскресенье)"
$ws.Range("V2").Value = "2024-08-19 (Понедельник)"
$ws.Range("W2").Value = "2024-08-20 (Вторник)"
$ws.Range("X2").Value = "2024-08-21 (Среда)"
$ws.Range("Y2").Value = "2024-08-22 (Четверг)"
$ws.Range("Z2").Value = "2024-08-23 (Пятница)"
$ws.Range("AA2").Value = "2024-08-24 (Суббота)"
$ws.Range("AB2").Value = "2024-08-25 (Воскресенье)"
$ws.Range("AC2").Value = "2024-08-26 (Понедельник)"
$ws.Range("AD2").Value = "2024-08-27 (Вторник)"
$ws.Range("AE2").Value = "2024-08-28 (Среда)"
$ws.Range("AF2").Value = "2024-08-29 (Четверг)"
$ws.Range("AG2").Value = "2024-08-30 (Пятница)"
$ws.Range("AH2").Value = "2024-08-31 (Суббота)"
$ws.Range("AI2").Value = "Переработка (ч)"

# -----------------------------------------------------------------
# 3. Existing employee's row-number label fix (2 -> 1)
# -----------------------------------------------------------------
$ws.Range("A3").Value = "1"

# -----------------------------------------------------------------
# 4. Existing employee's overtime total was a stray "100"; normalise
#    it to the same "hh:mm" formatting used by every other hours cell.
# -----------------------------------------------------------------
$ws.Range("AI3").Value = "00:00"

# -----------------------------------------------------------------
# 5. New employee block (rows 6:8): mirror the 3-row (Время прихода /
#    Время ухода / Кол-во часов) template used by rows 3:5, including
#    its merged cells and borders, then fill in ilyas's data.
#
#    Merge the destination cells *before* copying the template's
#    formatting across - merging an already-styled range recomputes
#    borders on the anchor cell (spawning a slightly different style),
#    while merging first and formatting after reuses the template's
#    existing style indices verbatim.
# -----------------------------------------------------------------
$ws.Range("A6:A8").MergeCells = $true
$ws.Range("B6:B8").MergeCells = $true
$ws.Range("AI6:AI8").MergeCells = $true

$ws.Range("A3:AI5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A6").Value = "2"
$ws.Range("B6").Value = "ilyas"
$ws.Range("C6").Value = "Время прихода"
$ws.Range("C7").Value = "Время ухода"
$ws.Range("C8").Value = "Кол-во часов"

$ws.Range("AB6").Value = "01:57"
$ws.Range("AI6").Value = "09:48"
$ws.Range("AB7").Value = "11:46"

$ws.Range("D8:AA8").Value = "00:00"
$ws.Range("AB8").Value = "09:49"
$ws.Range("AC8:AH8").Value = "00:00"
